$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E to be treated as plain text so that numeric-looking
# strings (prices like "592.79", "0.119", "7.80") are preserved exactly as
# text rather than being auto-converted to numbers by Excel.
$priceRange = $ws.Range("D2:E51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "63.500.04"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "3.176.57"
$ws.Range("E3").Value = "  -2.70%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "592.79"
$ws.Range("E5").Value = "  -1.25%  "
$ws.Range("D6").Value = "136.32"
$ws.Range("E6").Value = "  -1.35%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "3.172.89"
$ws.Range("E8").Value = "  -2.83%  "
$ws.Range("D9").Value = "0.512"
$ws.Range("E9").Value = "  -0.26%  "
$ws.Range("D10").Value = "0.143"
$ws.Range("E10").Value = "  -3.32%  "
$ws.Range("D11").Value = "5.34"
$ws.Range("E11").Value = "  -2.28%  "
$ws.Range("D12").Value = "0.456"
$ws.Range("E12").Value = "  -1.32%  "
$ws.Range("D13").Value = "0.0000240"
$ws.Range("E13").Value = "  -1.13%  "
$ws.Range("D14").Value = "34.73"
$ws.Range("E14").Value = "  +1.66%  "
$ws.Range("D15").Value = "3.696.27"
$ws.Range("E15").Value = "  -2.92%  "
$ws.Range("D16").Value = "0.119"
$ws.Range("E16").Value = "  -1.65%  "
$ws.Range("D17").Value = "3.175.78"
$ws.Range("E17").Value = "  -2.98%  "
$ws.Range("D18").Value = "63.431.10"
$ws.Range("E18").Value = "  +0.23%  "
$ws.Range("D19").Value = "6.54"
$ws.Range("E19").Value = "  -3.85%  "
$ws.Range("D20").Value = "461.13"
$ws.Range("E20").Value = "  -2.40%  "
$ws.Range("D21").Value = "13.93"
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("D22").Value = "0.695"
$ws.Range("E22").Value = "  -4.71%  "
$ws.Range("D23").Value = "7.65"
$ws.Range("E23").Value = "  -2.49%  "
$ws.Range("D24").Value = "83.07"
$ws.Range("E24").Value = "  -2.43%  "
$ws.Range("D25").Value = "13.19"
$ws.Range("E25").Value = "  -3.56%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").Value = "2.68"
$ws.Range("E27").Value = "  -2.32%  "
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.16%  "
$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").Value = "2.09"
$ws.Range("E29").Value = "  -1.38%  "
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").Value = "7.80"
$ws.Range("E30").Value = "  -3.17%  "
$ws.Range("D31").Value = "6.79"
$ws.Range("E31").Value = "  -2.86%  "
$ws.Range("D32").Value = "27.52"
$ws.Range("E32").Value = "  -2.71%  "
$ws.Range("D33").Value = "0.101"
$ws.Range("E33").Value = "  -2.08%  "
$ws.Range("D34").Value = "2.47"
$ws.Range("E34").Value = "  -0.44%  "
$ws.Range("D35").Value = "1.02"
$ws.Range("E35").Value = "  -5.09%  "
$ws.Range("D36").Value = "5.88"
$ws.Range("E36").Value = "  -1.16%  "
$ws.Range("D37").Value = "0.0₃0733"
$ws.Range("E37").Value = "  +0.93%  "
$ws.Range("D38").Value = "51.39"
$ws.Range("E38").Value = "  -1.29%  "
$ws.Range("D39").Value = "0.0392"
$ws.Range("E39").Value = "  -2.02%  "
$ws.Range("B40").Value = "Cosmos"
$ws.Range("C40").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D40").Value = "8.12"
$ws.Range("E40").Value = "  -1.28%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "2.67"
$ws.Range("E41").Value = "  -0.60%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "0.114"
$ws.Range("E42").Value = "  -4.36%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").Value = "396.32"
$ws.Range("E43").Value = "  -7.23%  "
$ws.Range("D44").Value = "2.807.09"
$ws.Range("E44").Value = "  -8.33%  "
$ws.Range("D45").Value = "0.253"
$ws.Range("E45").Value = "  -2.46%  "
$ws.Range("D46").Value = "129.36"
$ws.Range("E46").Value = "  +1.82%  "
$ws.Range("D47").Value = "2.14"
$ws.Range("E47").Value = "  -2.20%  "
$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").Value = "0.998"
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("B49").Value = "Arweave"
$ws.Range("C49").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D49").Value = "35.65"
$ws.Range("E49").Value = "  -1.68%  "
$ws.Range("D50").Value = "25.54"
$ws.Range("E50").Value = "  -2.11%  "
$ws.Range("D51").Value = "0.111"
$ws.Range("E51").Value = "  -1.84%  "

# Restore default (unstyled) formatting now that the text values are committed.
$priceRange.ClearFormats()

